$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds 3 repeated "group" blocks of columns (AUC, KS, KS_p_val,
# rel_type, GINI, Empty) across B:G, H:M, N:S. We are expanding each group
# to 9 columns (AUC, KS, KS_p_val, rel_type, GINI, Count, Empty,
# "Empty% in level", "Empty% in all Empty") by inserting 3 new columns into
# each group, right before its "Empty" column (this keeps the insertion
# strictly inside the row-2 group merge so Excel auto-extends it).
#
# Process the rightmost group first so earlier column letters stay stable.

$ws.Range("S:U").Insert()
$ws.Range("M:O").Insert()
$ws.Range("G:I").Insert()

# --- Row 3 header labels for each expanded group ---
$ws.Range("G3").Value = "Count"
$ws.Range("H3").Value = "Empty"
$ws.Range("I3").Value = "Empty% in level"
$ws.Range("J3").Value = "Empty% in all Empty"

$ws.Range("P3").Value = "Count"
$ws.Range("Q3").Value = "Empty"
$ws.Range("R3").Value = "Empty% in level"
$ws.Range("S3").Value = "Empty% in all Empty"

$ws.Range("Y3").Value = "Count"
$ws.Range("Z3").Value = "Empty"
$ws.Range("AA3").Value = "Empty% in level"
$ws.Range("AB3").Value = "Empty% in all Empty"

# --- Row 4 data values for each expanded group ---
# Group 1 (B:J)
$ws.Range("B4").Value = 0.5151886931344163
$ws.Range("C4").Value = 0.06024766485089394
$ws.Range("D4").Value = 0.3231387497223792
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.03037738626883257
$ws.Range("G4").Value = 574
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# Group 2 (K:S)
$ws.Range("K4").Value = 0.6916966180981214
$ws.Range("L4").Value = 0.2957859205901799
$ws.Range("M4").Value = [double]"4.412458392712792E-16"
$ws.Range("N4").Value = -1
$ws.Range("O4").Value = 0.3833932361962429
$ws.Range("P4").Value = 279
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0

# Group 3 (T:AB)
$ws.Range("T4").Value = 0.777910695344961
$ws.Range("U4").Value = 0.4515316091266518
$ws.Range("V4").Value = [double]"8.705849996169642E-24"
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 0.5558213906899221
$ws.Range("Y4").Value = 147
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0

# Trailing extra columns AC/AD (Empty / Empty part) row 4 values
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0

# --- Fix up the group "index" cells on row 2 (1 / 0 / 2) per group ---
$ws.Range("B2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("T2").Value = 2

Write-Host "done"
